function Set-CellText($ws, $ref, $val) {
    $ws.Range($ref).Value = "'" + $val
    $ws.Range($ref).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" '51.933.62'
Set-CellText $ws "E2" '  -14.60%  '
Set-CellText $ws "D3" '2.323.90'
Set-CellText $ws "E3" '  -20.06%  '
Set-CellText $ws "D4" '0.999'
Set-CellText $ws "E4" '  -0.10%  '
Set-CellText $ws "D5" '436.21'
Set-CellText $ws "E5" '  -17.07%  '
Set-CellText $ws "D6" '123.06'
Set-CellText $ws "E6" '  -14.61%  '
Set-CellText $ws "D7" '0.996'
Set-CellText $ws "E7" '  -0.23%  '
Set-CellText $ws "D8" '0.466'
Set-CellText $ws "E8" '  -14.60%  '
Set-CellText $ws "D9" '2.326.21'
Set-CellText $ws "E9" '  -20.18%  '
Set-CellText $ws "D10" '0.0896'
Set-CellText $ws "E10" '  -16.57%  '
Set-CellText $ws "D11" '5.14'
Set-CellText $ws "E11" '  -16.82%  '
Set-CellText $ws "D12" '0.306'
Set-CellText $ws "E12" '  -14.42%  '
Set-CellText $ws "E13" '  -6.34%  '
Set-CellText $ws "D14" '2.700.40'
Set-CellText $ws "E14" '  -20.93%  '
Set-CellText $ws "D15" '51.833.99'
Set-CellText $ws "E15" '  -14.78%  '
Set-CellText $ws "D16" '18.87'
Set-CellText $ws "E16" '  -16.31%  '
Set-CellText $ws "E17" '  -15.35%  '
Set-CellText $ws "D18" '2.312.69'
Set-CellText $ws "E18" '  -20.78%  '
Set-CellText $ws "D19" '4.06'
Set-CellText $ws "E19" '  -17.16%  '
Set-CellText $ws "D20" '296.62'
Set-CellText $ws "E20" '  -16.05%  '
Set-CellText $ws "D21" '8.90'
Set-CellText $ws "E21" '  -22.88%  '
Set-CellText $ws "D22" '1.01'
Set-CellText $ws "E22" '  +0.71%  '
Set-CellText $ws "E23" '  -0.18%  '
Set-CellText $ws "D24" '5.21'
Set-CellText $ws "E24" '  -19.96%  '
Set-CellText $ws "D25" '53.53'
Set-CellText $ws "E25" '  -17.45%  '
Set-CellText $ws "D26" '0.993'
Set-CellText $ws "E26" '  -0.67%  '
Set-CellText $ws "D27" '0.372'
Set-CellText $ws "E27" '  -17.22%  '
Set-CellText $ws "D28" '2.342.38'
Set-CellText $ws "E28" '  -22.78%  '
Set-CellText $ws "E29" '  -24.82%  '
Set-CellText $ws "B30" 'InternetComputer(DFINITY)'
Set-CellText $ws "C30" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws "D30" '6.86'
Set-CellText $ws "E30" '  -12.03%  '
Set-CellText $ws "B31" 'USDe'
Set-CellText $ws "C31" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-CellText $ws "D31" '0.996'
Set-CellText $ws "E31" '  -0.32%  '
Set-CellText $ws "D32" '0.0₃0667'
Set-CellText $ws "E32" '  -22.98%  '
Set-CellText $ws "D33" '143.37'
Set-CellText $ws "E33" '  -6.50%  '
Set-CellText $ws "D34" '17.07'
Set-CellText $ws "E34" '  -13.07%  '
Set-CellText $ws "E35" '  -20.78%  '
Set-CellText $ws "D36" '4.81'
Set-CellText $ws "E36" '  -13.86%  '
Set-CellText $ws "D37" '3.33'
Set-CellText $ws "E37" '  -24.19%  '
Set-CellText $ws "D38" '0.997'
Set-CellText $ws "E38" '  -16.74%  '
Set-CellText $ws "D39" '0.994'
Set-CellText $ws "E39" '  -0.32%  '
Set-CellText $ws "D40" '0.771'
Set-CellText $ws "E40" '  -22.44%  '
Set-CellText $ws "D41" '31.91'
Set-CellText $ws "E41" '  -15.11%  '
Set-CellText $ws "B42" 'Mantle'
Set-CellText $ws "C42" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText $ws "D42" '0.564'
Set-CellText $ws "E42" '  -13.57%  '
Set-CellText $ws "B43" 'WhiteBITCoin'
Set-CellText $ws "C43" 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-CellText $ws "D43" '10.11'
Set-CellText $ws "E43" '  -2.18%  '
Set-CellText $ws "B44" 'Filecoin'
Set-CellText $ws "C44" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws "D44" '3.18'
Set-CellText $ws "E44" '  -13.75%  '
Set-CellText $ws "D45" '0.0505'
Set-CellText $ws "E45" '  -13.20%  '
Set-CellText $ws "D46" '1.891.90'
Set-CellText $ws "E46" '  -17.41%  '
Set-CellText $ws "D47" '1.16'
Set-CellText $ws "E47" '  -21.15%  '
Set-CellText $ws "D48" '0.0205'
Set-CellText $ws "E48" '  -13.68%  '
Set-CellText $ws "D49" '0.0816'
Set-CellText $ws "E49" '  -10.54%  '
Set-CellText $ws "D50" '15.91'
Set-CellText $ws "E50" '  -21.91%  '
Set-CellText $ws "D51" '4.02'
Set-CellText $ws "E51" '  -18.57%  '
